# Update the "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to reflect the latest generated output (gh-pages update at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value  = 1151
$wsExpo.Range("F4").Value  = 260
$wsExpo.Range("F5").Value  = 144
$wsExpo.Range("F7").Value  = 12183
$wsExpo.Range("F8").Value  = 53
$wsExpo.Range("F10").Value = 118
$wsExpo.Range("F11").Value = 11967
$wsExpo.Range("F12").Value = 4790
$wsExpo.Range("F13").Value = 1790
$wsExpo.Range("F14").Value = 108
$wsExpo.Range("F15").Value = 42
$wsExpo.Range("F17").Value = 90
$wsExpo.Range("F18").Value = 941
$wsExpo.Range("F20").Value = 161
$wsExpo.Range("F21").Value = 67

# --- Sheet "全部类型" (All types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value  = 1151
$wsAll.Range("F4").Value  = 260
$wsAll.Range("F5").Value  = 144
$wsAll.Range("F9").Value  = 12183
$wsAll.Range("F10").Value = 53
$wsAll.Range("F12").Value = 118
$wsAll.Range("F13").Value = 11967
$wsAll.Range("F14").Value = 4790
$wsAll.Range("F15").Value = 1791
$wsAll.Range("F16").Value = 108
$wsAll.Range("F17").Value = 42
$wsAll.Range("F19").Value = 90
$wsAll.Range("F20").Value = 941
$wsAll.Range("F22").Value = 161
$wsAll.Range("F23").Value = 67
